$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (9-16) to append below the existing table.
# Columns: A,B,C,D,E,F,G,H,I,J,K,L,M,N
$rows = @(
    @{ r=9;  A=20200107; B=14; C=5; D=2; K=-961; L=1282.28;            M=858.35;  N=46 },
    @{ r=10; A=20200107; B=14; C=5; D=2; K=-960; L=1282.9100000000001; M=937.86;  N=45.87 },
    @{ r=11; A=20200107; B=14; C=5; D=2; K=-959; L=1282.0999999999999; M=1018.4;  N=44.41 },
    @{ r=12; A=20200107; B=14; C=5; D=2; E=2020; F=3; G=10; H=17; I=21; J=9.94; K=-911; L=1285; M=1063 },
    @{ r=13; A=20200107; B=14; C=5; D=2; E=2020; F=3; G=10; H=17; I=21; J=9.94; K=-910; L=1265; M=1046 },
    @{ r=14; A=20200107; B=14; C=5; D=2; E=2020; F=3; G=10; H=17; I=21; J=9.94; K=-909; L=1261; M=1026 },
    @{ r=15; A=20200107; B=14; C=5; D=2; E=2020; F=3; G=10; H=17; I=21; J=9.94; K=-908; L=1290; M=1006 },
    @{ r=16; A=20200107; B=14; C=5; D=2; E=2020; F=3; G=10; H=17; I=21; J=9.94; K=-907; L=1282; M=982 }
)

foreach ($row in $rows) {
    $r = $row.r
    if ($row.ContainsKey('A')) { $ws.Cells.Item($r, 1).Value = $row.A }
    if ($row.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($row.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $row.C }
    if ($row.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = $row.D }
    if ($row.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $row.E }
    if ($row.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $row.F }
    if ($row.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = $row.G }
    if ($row.ContainsKey('H')) { $ws.Cells.Item($r, 8).Value = $row.H }
    if ($row.ContainsKey('I')) { $ws.Cells.Item($r, 9).Value = $row.I }
    if ($row.ContainsKey('J')) {
        $cell = $ws.Cells.Item($r, 10)
        $cell.Value = $row.J
        $cell.NumberFormat = "0.00E+00"
    }
    if ($row.ContainsKey('K')) { $ws.Cells.Item($r, 11).Value = $row.K }
    if ($row.ContainsKey('L')) { $ws.Cells.Item($r, 12).Value = $row.L }
    if ($row.ContainsKey('M')) { $ws.Cells.Item($r, 13).Value = $row.M }
    if ($row.ContainsKey('N')) { $ws.Cells.Item($r, 14).Value = $row.N }
}

# Update the selected cell to match the diff's new selection.
$ws.Range("M22").Select()
